$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 3, pushing existing rows 3-8 down to 5-10
$ws.Range("A3:A4").EntireRow.Insert()

# Row 3: new combined parallel testing entry
$ws.Cells.Item(3, 1).Value = "Combined Youden's J Optimized with CXR (Parallel)"
$ws.Cells.Item(3, 2).Value = 94
$ws.Cells.Item(3, 3).Value = 76.5

# Row 4: new serial testing entry
$ws.Cells.Item(4, 1).Value = "VOC Model (TPP Optimized)"
$ws.Cells.Item(4, 2).Value = 84.40000000000001
$ws.Cells.Item(4, 3).Value = 66.7

# Row 5: rename only (values retained from the old "VOC Model (TPP Threshold)" row)
$ws.Cells.Item(5, 1).Value = "VOC Model (Sensitivity Optimized)"

# Row 6: rename only (values retained from the old "VOC Model (Optimized)" row)
$ws.Cells.Item(6, 1).Value = "VOC Model (Youden's J)"

# Rows 7-10 (CXR + CAD, CXR, CRP, Urine LAM) keep their original labels/values;
# they were already shifted into place by the Insert() above.
